$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set cells in the same order the original author likely typed them, so that
# new shared-string entries are appended in the same order as the target
# workbook (L, K, D, J are new strings; B, C, E reuse existing strings).
$ws.Range("L77").Value = "logs-hjxh-2018-4-27-class14-original-more-percent64"
$ws.Range("K77").Value = "python feed_run.py --output_mode=classes --output_nodes=14 --input_nums=92 --input_nodes=92 --low_nums=2 --low_nodes=92 --low_fun=elu --one_hot=True --input_fun=tanh --batch_size=100 --learning_rate=0.001 --train_mode=Adadelta --eval_size=5400 --test_size=1339 --use_biases=yes --dropout_in=0.6 --dropout_low=0.6"
$ws.Range("D77").Value = "batch_size=100 low_nums=2 use_biases=yes  dropout_in=0.6  dropout_low=0.6"
$ws.Range("J77").Value = "经过约41小时，拟合精度和泛化精度接近，运行较久，但是没多少变化和进展。"

$ws.Range("A77").Value = 43217.356944444444
$ws.Range("B77").Value = "分类O"
$ws.Range("C77").Value = "14分类"
$ws.Range("E77").Value = "最高标签，重新训练，原始数据新加指标train-hjxh365-2018-4-16-day-high-original-more"
$ws.Range("F77").Value = 0.63
$ws.Range("G77").Value = 0.64
$ws.Range("H77").Value = 0.95
$ws.Range("I77").Value = 0.94

$ws.Rows.Item(77).RowHeight = 82.5

$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Application.ActiveWindow.ScrollColumn = 4

Write-Output "done"
